# Update automatic: dades i banners [2026-02-16 20:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "'2026-02-16 20:18:47"
$ws.Range("H2").Value = "'96%"
$ws.Range("I2").Value = "'20.6 mm"
$ws.Range("E3").Value = "'2026-02-16 20:18:49"
$ws.Range("G3").Value = "'237 cm"
$ws.Range("I3").Value = "'10.0 mm"
$ws.Range("N3").Value = "'-2.2 °C 19:49 TU"
$ws.Range("E4").Value = "'2026-02-16 20:18:52"
$ws.Range("H4").Value = "'60%"
$ws.Range("K4").Value = "'13.4 MJ/m2"
$ws.Range("O4").Value = "'13.8 °C"
$ws.Range("E5").Value = "'2026-02-16 20:18:54"
$ws.Range("I5").Value = "'24.6 mm"
$ws.Range("N5").Value = "'-1.9 °C 19:57 TU"
$ws.Range("O5").Value = "'-0.7 °C"
$ws.Range("E6").Value = "'2026-02-16 20:18:57"
$ws.Range("E7").Value = "'2026-02-16 20:18:59"
$ws.Range("J7").Value = "'1013.3 hPa"
$ws.Range("E8").Value = "'2026-02-16 20:19:02"
$ws.Range("K8").Value = "'12.0 MJ/m2"
$ws.Range("E9").Value = "'2026-02-16 20:19:04"
$ws.Range("O9").Value = "'11.4 °C"
$ws.Range("E10").Value = "'2026-02-16 20:19:07"
$ws.Range("E11").Value = "'2026-02-16 20:19:10"
$ws.Range("O11").Value = "'6.9 °C"
$ws.Range("E12").Value = "'2026-02-16 20:19:12"
$ws.Range("H12").Value = "'79%"
$ws.Range("E13").Value = "'2026-02-16 20:19:14"
$ws.Range("J13").Value = "'1014.7 hPa"
$ws.Range("O13").Value = "'5.8 °C"
$ws.Range("E14").Value = "'2026-02-16 20:19:17"
$ws.Range("O14").Value = "'16.2 °C"
$ws.Range("E15").Value = "'2026-02-16 20:19:19"
$ws.Range("O15").Value = "'11.6 °C"
$ws.Range("E16").Value = "'2026-02-16 20:19:22"
$ws.Range("N16").Value = "'-2.1 °C 19:49 TU"
$ws.Range("E17").Value = "'2026-02-16 20:19:24"
$ws.Range("E18").Value = "'2026-02-16 20:19:27"
$ws.Range("E19").Value = "'2026-02-16 20:19:29"
$ws.Range("O19").Value = "'7.1 °C"
$ws.Range("E20").Value = "'2026-02-16 20:19:32"
$ws.Range("I20").Value = "'0.5 mm"
$ws.Range("E21").Value = "'2026-02-16 20:19:34"
$ws.Range("H21").Value = "'70%"
$ws.Range("J21").Value = "'1014.2 hPa"
$ws.Range("O21").Value = "'8.6 °C"
$ws.Range("E22").Value = "'2026-02-16 20:19:37"
$ws.Range("E23").Value = "'2026-02-16 20:19:39"
$ws.Range("I23").Value = "'14.4 mm"
$ws.Range("N23").Value = "'-2.6 °C 19:59 TU"
$ws.Range("E24").Value = "'2026-02-16 20:19:42"
$ws.Range("J24").Value = "'1016.5 hPa"
$ws.Range("E25").Value = "'2026-02-16 20:19:45"
$ws.Range("I25").Value = "'6.1 mm"
$ws.Range("N25").Value = "'-0.7 °C 19:31 TU"
$ws.Range("E26").Value = "'2026-02-16 20:19:47"
$ws.Range("E27").Value = "'2026-02-16 20:19:50"
$ws.Range("O27").Value = "'1.2 °C"
$ws.Range("E28").Value = "'2026-02-16 20:19:52"
$ws.Range("K28").Value = "'12.8 MJ/m2"
$ws.Range("E29").Value = "'2026-02-16 20:19:55"
$ws.Range("H29").Value = "'80%"
$ws.Range("E30").Value = "'2026-02-16 20:19:57"
$ws.Range("E31").Value = "'2026-02-16 20:20:00"
$ws.Range("E32").Value = "'2026-02-16 20:20:02"
$ws.Range("H32").Value = "'81%"
$ws.Range("E33").Value = "'2026-02-16 20:20:05"
$ws.Range("E34").Value = "'2026-02-16 20:20:07"
$ws.Range("E35").Value = "'2026-02-16 20:20:10"
$ws.Range("H35").Value = "'73%"
$ws.Range("I35").Value = "'0.2 mm"
$ws.Range("O35").Value = "'9.5 °C"
$ws.Range("E36").Value = "'2026-02-16 20:20:13"
$ws.Range("H36").Value = "'73%"
$ws.Range("L36").Value = "'62.3 km/h - 313º 19:54 TU"
$ws.Range("O36").Value = "'11.9 °C"
$ws.Range("E37").Value = "'2026-02-16 20:20:15"
$ws.Range("J37").Value = "'1014.7 hPa"
$ws.Range("E38").Value = "'2026-02-16 20:20:18"
$ws.Range("E39").Value = "'2026-02-16 20:20:20"
$ws.Range("I39").Value = "'3.8 mm"
$ws.Range("N39").Value = "'-1.8 °C 19:59 TU"
$ws.Range("O39").Value = "'0.3 °C"
$ws.Range("E40").Value = "'2026-02-16 20:20:23"
$ws.Range("J40").Value = "'1016.3 hPa"
$ws.Range("E41").Value = "'2026-02-16 20:20:25"
$ws.Range("J41").Value = "'1014.8 hPa"
$ws.Range("E42").Value = "'2026-02-16 20:20:28"
$ws.Range("H42").Value = "'80%"
$ws.Range("E43").Value = "'2026-02-16 20:20:30"
$ws.Range("O43").Value = "'8.8 °C"
$ws.Range("E44").Value = "'2026-02-16 20:20:33"
$ws.Range("G44").Value = "'243 cm"
$ws.Range("I44").Value = "'10.1 mm"
$ws.Range("L44").Value = "'64.1 km/h - 75º 19:44 TU"
$ws.Range("N44").Value = "'-2.1 °C 19:58 TU"
$ws.Range("E45").Value = "'2026-02-16 20:20:35"
$ws.Range("I45").Value = "'17.7 mm"
$ws.Range("E46").Value = "'2026-02-16 20:20:38"
$ws.Range("O46").Value = "'16.2 °C"
